$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.944.90"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.845.05"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "702.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").Value = "3.843.98"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "4.490.02"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "3.840.25"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "70.918.45"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.182"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("D35").Value = "3.798.60"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.96%  "
$ws.Range("E40").Value = "  +6.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000314"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "406.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.36%  "
